$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.717.81"
$ws.Range("E2").Value = "  +2.84%  "

$ws.Range("D3").Value = "'2.525.99"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'592.87"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'175.66"
$ws.Range("E6").Value = "  +3.71%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").Value = "'2.525.99"
$ws.Range("E9").Value = "  +1.74%  "

$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  +3.26%  "

$ws.Range("E11").Value = "  +3.20%  "

$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "'26.85"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "'2.982.59"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").Value = "'0.0000178"
$ws.Range("E16").Value = "  +3.06%  "

$ws.Range("D17").Value = "'67.652.70"
$ws.Range("E17").Value = "  +3.06%  "

$ws.Range("D18").Value = "'2.510.68"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").Value = "'7.99"
$ws.Range("E19").Value = "  +5.88%  "

$ws.Range("D20").Value = "'11.43"
$ws.Range("E20").Value = "  +2.96%  "

$ws.Range("D21").Value = "'362.21"
$ws.Range("E21").Value = "  +5.80%  "

$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "'4.66"
$ws.Range("E23").Value = "  +3.19%  "

$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "'71.07"
$ws.Range("E26").Value = "  +3.61%  "

$ws.Range("D27").Value = "'10.25"
$ws.Range("E27").Value = "  +3.96%  "

$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").Value = "'2.652.70"
$ws.Range("E29").Value = "  +1.46%  "

$ws.Range("D30").Value = "'0.0₃0990"
$ws.Range("E30").Value = "  +3.59%  "

$ws.Range("D31").Value = "'542.82"
$ws.Range("E31").Value = "  +4.58%  "

$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "  +3.51%  "

$ws.Range("E33").Value = "  +2.86%  "

$ws.Range("E34").Value = "  +3.52%  "

$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "'1.46"
$ws.Range("E37").Value = "  +1.65%  "

$ws.Range("D38").Value = "'156.49"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("D39").Value = "'18.83"
$ws.Range("E39").Value = "  +2.35%  "

$ws.Range("D40").Value = "'18.66"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").Value = "'5.18"
$ws.Range("E42").Value = "  +3.25%  "

$ws.Range("E43").Value = "  +2.41%  "

$ws.Range("E44").Value = "  +4.54%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").Value = "'146.44"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("E48").Value = "  +2.38%  "

$ws.Range("D49").Value = "'0.0₆0277"
$ws.Range("E49").Value = "  +3.62%  "

$ws.Range("D50").Value = "'1.69"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("E51").Value = "  +1.73%  "
